$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "vessels" sheet: add new "cost_operation" column J
# ---------------------------------------------------------------------------
$wsVessels = $wb.Worksheets.Item("vessels")

$wsVessels.Range("J1").Value = "cost_operation"
$wsVessels.Range("J2").Value = 100
$wsVessels.Range("J3").Value = 100
$wsVessels.Range("J4").Value = 100

# best-effort column widths
$wsVessels.Columns.Item(6).ColumnWidth = 11.6640625
$wsVessels.Columns.Item(8).ColumnWidth = 16.21875
$wsVessels.Columns.Item(10).ColumnWidth = 14.5546875

# ---------------------------------------------------------------------------
# 2. "general" sheet: add new parameter columns D:H with headers + values
# ---------------------------------------------------------------------------
$wsGeneral = $wb.Worksheets.Item("general")

$wsGeneral.Range("D1").Value = "cost_technicians"
$wsGeneral.Range("E1").Value = "cost_downtime"
$wsGeneral.Range("F1").Value = "penalty_cost_late"
$wsGeneral.Range("G1").Value = "penalty_cost_not_performed"
$wsGeneral.Range("H1").Value = "latest_period"

$wsGeneral.Range("D2").Value = 50
$wsGeneral.Range("E2").Value = 500
$wsGeneral.Range("F2").Value = 20000
$wsGeneral.Range("G2").Value = 1000000
$wsGeneral.Range("H2").Value = 80

# best-effort column widths (matches the bestFit widths Excel computed)
$wsGeneral.Columns.Item(4).ColumnWidth = 14.6640625
$wsGeneral.Columns.Item(5).ColumnWidth = 13.21875
$wsGeneral.Columns.Item(6).ColumnWidth = 20.21875
$wsGeneral.Columns.Item(7).ColumnWidth = 25.6640625
$wsGeneral.Columns.Item(8).ColumnWidth = 11.33203125

# page setup (A4, portrait)
$wsGeneral.PageSetup.PaperSize = 9
$wsGeneral.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 3. "tasks" sheet: column widths only (no content change)
# ---------------------------------------------------------------------------
$wsTasks = $wb.Worksheets.Item("tasks")

$wsTasks.Columns.Item(3).ColumnWidth = 10.109375
$wsTasks.Columns.Item(4).ColumnWidth = 10.33203125
$wsTasks.Columns.Item(5).ColumnWidth = 10.109375
$wsTasks.Columns.Item(6).ColumnWidth = 10

# ---------------------------------------------------------------------------
# 4. Selections on each sheet (so the saved file keeps the right
#    activeCell / sqref per sheet), applied in the same order the user
#    would have tabbed through them.
# ---------------------------------------------------------------------------

# bases: selection moved to C1
$wsBases = $wb.Worksheets.Item("bases")
[void]$wsBases.Select()
[void]$wsBases.Range("C1").Select()

# vessels: whole column F selected
[void]$wsVessels.Select()
[void]$wsVessels.Columns.Item(6).Select()

# tasks: F2 selected
[void]$wsTasks.Select()
[void]$wsTasks.Range("F2").Select()

# general: H3 selected, and ends up the active sheet/tab
[void]$wsGeneral.Select()
[void]$wsGeneral.Range("H3").Select()
